$d = $word.ActiveDocument

$pairs = @(
    @("92×44=4048", "45×20=900"),
    @("20×32=640", "89×46=4094"),
    @("65×22=1430", "68×17=1156"),
    @("50×79=3950", "32×82=2624"),
    @("57×43=2451", "68×90=6120"),
    @("26×54=1404", "29×98=2842"),
    @("64×64=4096", "50×71=3550"),
    @("61×71=4331", "27×69=1863"),
    @("97×72=6984", "78×65=5070"),
    @("38×97=3686", "97×12=1164"),
    @("53×51=2703", "11×73=803"),
    @("76×45=3420", "27×84=2268"),
    @("67×71=4757", "63×49=3087"),
    @("92×27=2484", "84×95=7980"),
    @("45×73=3285", "31×95=2945"),
    @("78×73=5694", "70×90=6300"),
    @("20×17=340", "74×19=1406"),
    @("63×98=6174", "96×14=1344"),
    @("32×79=2528", "95×30=2850"),
    @("29×93=2697", "72×51=3672"),
    @("79×71=5609", "60×41=2460"),
    @("63×61=3843", "37×12=444"),
    @("48×73=3504", "83×22=1826"),
    @("43×85=3655", "82×68=5576"),
    @("75×96=7200", "59×53=3127")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
